$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.216.31"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").Value = "3.997.38"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.94"
$ws.Range("E5").Value = "  +14.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.82"
$ws.Range("E6").Value = "  +11.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.685"
$ws.Range("E7").Value = "  -1.85%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +1.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.93"
$ws.Range("E11").Value = "  +6.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000339"
$ws.Range("E12").Value = "  +3.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.13"
$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").Value = "4.624.35"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").Value = "3.986.58"
$ws.Range("E15").Value = "  -1.16%  "

$ws.Range("E16").Value = "  +4.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.27"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.59"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "73.037.70"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.46"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.94"
$ws.Range("E22").Value = "  +15.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.24"
$ws.Range("E23").Value = "  -2.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.38"
$ws.Range("E24").Value = "  -4.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.23"
$ws.Range("E25").Value = "  -2.98%  "

$ws.Range("E26").Value = "  -6.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.13"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.59"
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.95"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.24"
$ws.Range("E30").Value = "  -2.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.71"
$ws.Range("E31").Value = "  -7.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.75"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "72.41"
$ws.Range("E34").Value = "  +7.61%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0000103"
$ws.Range("E35").Value = "  +19.17%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.09"
$ws.Range("E36").Value = "  -4.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "634.56"
$ws.Range("E37").Value = "  -6.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.432"
$ws.Range("E38").Value = "  -6.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.146"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.08"
$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  -4.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0486"
$ws.Range("E45").Value = "  -1.70%  "

$ws.Range("E46").Value = "  -1.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  +2.34%  "

$ws.Range("E48").Value = "  -2.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  +31.21%  "

$ws.Range("D50").Value = "2.887.04"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("E51").Value = "  -2.27%  "

